$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update claim number value in B2 (text value, keep as text since it's numeric-looking)
$ws.Range("B2").Value = "'1120194100385"

# Update selection to F6
$ws.Range("F6").Select()
